$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Mapping of row -> (D new value, E new value). $null means "leave unchanged".
$updates = @{
    2  = @("299.14", "-0.98%")
    3  = @("31.44", "-0.16%")
    4  = @("5.115", "-0.40%")
    5  = @("0.08004", "8.27%")
    6  = @("2.395", "8.60%")
    7  = @("7.781", "-1.53%")
    8  = @("3.907", "2.48%")
    9  = @("0.9299", "1.08%")
    10 = @("0.1741", "0.84%")
    11 = @("0.07460", "0.16%")
    12 = @("0.09118", "11.28%")
    13 = @("0.03025", "1.24%")
    14 = @("0.1001", "0.96%")
    15 = @("0.001512", "1.28%")
    16 = @("0.006198", "1.09%")
    17 = @("3.497", "1.07%")
    18 = @("2.292", "2.86%")
    19 = @($null, "0.30%")
    20 = @("0.1332", "-0.32%")
    21 = @("4.161", "-10.26%")
    22 = @("0.1714", "9.32%")
    23 = @("0.04631", "0.28%")
    24 = @("0.001245", "1.78%")
    25 = @("0.004481", "0.04%")
    26 = @("0.0001201", "-7.32%")
    27 = @("0.0003418", "24.87%")
    28 = @($null, $null)
    29 = @($null, $null)
    30 = @($null, $null)
    31 = @($null, $null)
    32 = @($null, $null)
    33 = @($null, $null)
    34 = @($null, $null)
    35 = @($null, $null)
    36 = @($null, $null)
    37 = @($null, $null)
    38 = @($null, $null)
    39 = @("0.01751", "1.72%")
    40 = @("0.04589", "1.51%")
    41 = @("0.006914", "-5.54%")
    42 = @("0.1362", "0.94%")
    43 = @("0.002144", "0.00%")
    44 = @("0.01039", "-2.41%")
    45 = @("0.00006310", "0.76%")
    46 = @("0.00000000752", "0.22%")
    47 = @("0.008018", "-19.75%")
    48 = @("0.7482", "-8.83%")
    49 = @("0.00002105", "0.22%")
    50 = @("0.0002005", "0.29%")
    51 = @($null, $null)
}

for ($row = 2; $row -le 51; $row++) {
    $pair = $updates[$row]
    $dVal = $pair[0]
    $eVal = $pair[1]

    if ($null -ne $dVal) {
        $cell = $ws.Cells.Item($row, 4)
        $cell.NumberFormat = "@"
        $cell.Value = $dVal
    }
    if ($null -ne $eVal) {
        $cell = $ws.Cells.Item($row, 5)
        $cell.NumberFormat = "@"
        $cell.Value = $eVal
    }

    $cellF = $ws.Cells.Item($row, 6)
    $cellF.NumberFormat = "@"
    $cellF.Value = "17-1-2023"

    $cellG = $ws.Cells.Item($row, 7)
    $cellG.NumberFormat = "@"
    $cellG.Value = "0"
}
